$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Kitchen & Dining / Kitchen & Dining -> frequency 964 -> 936
$ws.Range("C3").Value = 936

# Row 4: Storage & Organization -> Fashion & Accessories ; 522 -> 508
$ws.Range("A4").Value = "Fashion & Accessories"
$ws.Range("B4").Value = "Fashion & Accessories"
$ws.Range("C4").Value = 508

# Row 5: Fashion & Accessories -> Storage & Organization ; 508 -> 490
$ws.Range("A5").Value = "Storage & Organization"
$ws.Range("B5").Value = "Storage & Organization"
$ws.Range("C5").Value = 490

# Row 6: Home Decor / Home Decor -> frequency 410 -> 400
$ws.Range("C6").Value = 400

# Row 7: Kids & Toys / Storage & Organization -> Storage & Organization / Kids & Toys ; 182 -> 158
$ws.Range("A7").Value = "Storage & Organization"
$ws.Range("B7").Value = "Kids & Toys"
$ws.Range("C7").Value = 158

# Row 8: Storage & Organization / Kids & Toys -> Kids & Toys / Storage & Organization ; 182 -> 158
$ws.Range("A8").Value = "Kids & Toys"
$ws.Range("B8").Value = "Storage & Organization"
$ws.Range("C8").Value = 158

# Row 9: Storage & Organization / Fashion & Accessories -> frequency 156 -> 150
$ws.Range("C9").Value = 150

# Row 10: Fashion & Accessories / Storage & Organization -> frequency 156 -> 150
$ws.Range("C10").Value = 150

# Row 11: Kitchen & Dining / Storage & Organization -> Kids & Toys / Kids & Toys ; 128 -> 120
$ws.Range("A11").Value = "Kids & Toys"
$ws.Range("B11").Value = "Kids & Toys"
$ws.Range("C11").Value = 120
